# Trade #118 closed at 2026-02-18 00:38:50 - unknown UNKNOWN +0.000%
#
# This script updates the "live_trading_results" workbook:
#  - Refreshes the Summary sheet aggregate metrics
#  - Refreshes the Strategy Status row for HighProbConvergence
#  - Marks the open HighProbConvergence trade (#146) as CLOSED with an
#    early_exit reason, on both the "All Trades" sheet and its
#    strategy-specific "HighProbConvergence" sheet
#  - Appends two newly opened trades (#175 momentum, #176 MarketMaking)
#    to "All Trades" and to their respective strategy sheets

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1499.02
$summary.Range("B4").Value = 0.13
$summary.Range("B5").Value = 0.02
$summary.Range("B6").Value = 146
$summary.Range("B8").Value = 51
$summary.Range("B9").Value = 45.89

# ---------------------------------------------------------------------
# Strategy Status sheet - row 3 (HighProbConvergence)
# ---------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C3").Value = 100.3
$status.Range("D3").Value = 18
$status.Range("E3").Value = 0.31
$status.Range("F3").Value = 0.3
$status.Range("G3").Value = 61.11

# ---------------------------------------------------------------------
# All Trades sheet - close out trade #146 (row 147)
# ---------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
$allTrades.Cells.Item(147, 7).Value = 0.044985
$allTrades.Cells.Item(147, 8).Value = "CLOSED"
$allTrades.Cells.Item(147, 9).Value = -65.11839999999999
$allTrades.Cells.Item(147, 10).Value = -0.08
$allTrades.Cells.Item(147, 11).Value = 100.3
$allTrades.Cells.Item(147, 12).Value = "early_exit"
$allTrades.Cells.Item(147, 13).Value = 0.17

# ---------------------------------------------------------------------
# HighProbConvergence sheet - close out the same trade (row 19)
# ---------------------------------------------------------------------
$hpc = $wb.Worksheets.Item("HighProbConvergence")
$hpc.Cells.Item(19, 7).Value = 0.044985
$hpc.Cells.Item(19, 8).Value = "CLOSED"
$hpc.Cells.Item(19, 9).Value = -65.11839999999999
$hpc.Cells.Item(19, 10).Value = -0.08
$hpc.Cells.Item(19, 11).Value = 100.3
$hpc.Cells.Item(19, 16).Value = "early_exit"
$hpc.Cells.Item(19, 17).Value = 0.17

# ---------------------------------------------------------------------
# Helper: write a text value while forcing text format so Excel does
# not auto-convert date/time-looking strings into date/time serials.
# ---------------------------------------------------------------------
function Set-TextCell($cell, $text) {
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# All Trades sheet - append new trade #175 (momentum) as row 176
# ---------------------------------------------------------------------
$row = $allTrades.Cells.Item(176, 1)
$row.Value = 175
Set-TextCell $allTrades.Cells.Item(176, 2) "2026-02-18"
Set-TextCell $allTrades.Cells.Item(176, 3) "00:38:43"
Set-TextCell $allTrades.Cells.Item(176, 4) "momentum"
Set-TextCell $allTrades.Cells.Item(176, 5) "DOWN"
$allTrades.Cells.Item(176, 6).Value = 0.128966
Set-TextCell $allTrades.Cells.Item(176, 8) "OPEN"
$allTrades.Cells.Item(176, 9).Value = 0
$allTrades.Cells.Item(176, 10).Value = 0
$allTrades.Cells.Item(176, 11).Value = 99.22374292899114
$allTrades.Cells.Item(176, 13).Value = 0
$allTrades.Cells.Item(176, 14).Value = 0
$allTrades.Cells.Item(176, 15).Value = 0
$allTrades.Cells.Item(176, 16).Value = 0.9
Set-TextCell $allTrades.Cells.Item(176, 17) "Downward momentum: -1.942% over 10 samples"

# ---------------------------------------------------------------------
# All Trades sheet - append new trade #176 (MarketMaking) as row 177
# ---------------------------------------------------------------------
$allTrades.Cells.Item(177, 1).Value = 176
Set-TextCell $allTrades.Cells.Item(177, 2) "2026-02-18"
Set-TextCell $allTrades.Cells.Item(177, 3) "00:38:44"
Set-TextCell $allTrades.Cells.Item(177, 4) "MarketMaking"
Set-TextCell $allTrades.Cells.Item(177, 5) "UP"
$allTrades.Cells.Item(177, 6).Value = 0.88
Set-TextCell $allTrades.Cells.Item(177, 8) "OPEN"
$allTrades.Cells.Item(177, 9).Value = 0
$allTrades.Cells.Item(177, 10).Value = 0
$allTrades.Cells.Item(177, 11).Value = 99.21858346467945
$allTrades.Cells.Item(177, 13).Value = 0
$allTrades.Cells.Item(177, 14).Value = 0
$allTrades.Cells.Item(177, 15).Value = 0
$allTrades.Cells.Item(177, 16).Value = 0.6
Set-TextCell $allTrades.Cells.Item(177, 17) "Normal spread capture: 198 bps"

# ---------------------------------------------------------------------
# momentum sheet - append new trade #175 as row 47
# ---------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")
$momentum.Cells.Item(47, 1).Value = 175
Set-TextCell $momentum.Cells.Item(47, 2) "2026-02-18"
Set-TextCell $momentum.Cells.Item(47, 3) "00:38:43"
Set-TextCell $momentum.Cells.Item(47, 4) "momentum"
Set-TextCell $momentum.Cells.Item(47, 5) "DOWN"
$momentum.Cells.Item(47, 6).Value = 0.128966
Set-TextCell $momentum.Cells.Item(47, 8) "OPEN"
$momentum.Cells.Item(47, 9).Value = 0
$momentum.Cells.Item(47, 10).Value = 0
$momentum.Cells.Item(47, 11).Value = 99.22374292899114
$momentum.Cells.Item(47, 12).Value = 0
$momentum.Cells.Item(47, 13).Value = 0
$momentum.Cells.Item(47, 14).Value = 0.9
Set-TextCell $momentum.Cells.Item(47, 15) "Downward momentum: -1.942% over 10 samples"
$momentum.Cells.Item(47, 17).Value = 0

# ---------------------------------------------------------------------
# MarketMaking sheet - append new trade #176 as row 70
# ---------------------------------------------------------------------
$mm = $wb.Worksheets.Item("MarketMaking")
$mm.Cells.Item(70, 1).Value = 176
Set-TextCell $mm.Cells.Item(70, 2) "2026-02-18"
Set-TextCell $mm.Cells.Item(70, 3) "00:38:44"
Set-TextCell $mm.Cells.Item(70, 4) "MarketMaking"
Set-TextCell $mm.Cells.Item(70, 5) "UP"
$mm.Cells.Item(70, 6).Value = 0.88
Set-TextCell $mm.Cells.Item(70, 8) "OPEN"
$mm.Cells.Item(70, 9).Value = 0
$mm.Cells.Item(70, 10).Value = 0
$mm.Cells.Item(70, 11).Value = 99.21858346467945
$mm.Cells.Item(70, 12).Value = 0
$mm.Cells.Item(70, 13).Value = 0
$mm.Cells.Item(70, 14).Value = 0.6
Set-TextCell $mm.Cells.Item(70, 15) "Normal spread capture: 198 bps"
$mm.Cells.Item(70, 17).Value = 0

$wb.Save()
